$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "Q2" = 1.98
    "R2" = 1.88
    "AN2" = 3.2
    "AZ2" = 201
    "Q3" = 2.3
    "R3" = 1.6
    "J5" = 2.25
    "O5" = 1.3
    "P5" = 3.4
    "R5" = 1.8
    "X5" = 7.5
    "AC5" = 9.5
    "AF5" = 51
    "AL5" = 41
    "AN5" = 3.6
    "G11" = 4.35
    "H11" = 3.25
    "I11" = 1.82
    "J11" = 4.75
    "K11" = 2.02
    "L11" = 2.4
    "M11" = 1.03
    "N11" = 6.65
    "O11" = 1.4
    "P11" = 2.5
    "Q11" = 2.18
    "U11" = 2.02
    "V11" = 1.62
    "X11" = 23
    "Y11" = 15
    "Z11" = 75
    "AA11" = 50
    "AB11" = 65
    "AC11" = 7.4
    "AD11" = 6.4
    "AE11" = 18.5
    "AG11" = 5.5
    "AH11" = 7.5
    "AJ11" = 14.5
    "AK11" = 17
    "AN11" = 5.9
    "AO11" = 26
    "AP11" = 35
    "AR11" = 200
    "AS11" = 500
    "AU11" = 7.8
    "AW11" = 3.5
    "AX11" = 9
    "AZ11" = 35
    "BA11" = 75
    "G14" = 2.25
    "I14" = 3.2
    "J14" = 3.2
    "K14" = 1.87
    "L14" = 4.33
    "M14" = 1.11
    "N14" = 6.5
    "O14" = 1.53
    "P14" = 2.38
    "W14" = 5.5
    "X14" = 9.5
    "Y14" = 10
    "Z14" = 21
    "AA14" = 23
    "AH14" = 15
    "AI14" = 13
    "AJ14" = 41
    "AL14" = 51
    "AN14" = 4
    "AP14" = 29
    "AR14" = 81
    "AS14" = 301
    "AW14" = 5
    "AX14" = 21
    "AZ14" = 81
    "G15" = 1.75
    "H15" = 3.2
    "I15" = 4.75
    "J15" = 2.5
    "L15" = 5.5
    "U15" = 2.2
    "V15" = 1.62
    "AD15" = 6.5
    "AG15" = 10
    "AI15" = 17
    "AK15" = 41
    "AN15" = 3.6
    "AO15" = 10
    "AZ15" = 101
    "G24" = 2.15
    "H24" = 3.1
    "I24" = 3.7
    "J24" = 2.88
    "K24" = 2.05
    "L24" = 4
    "X24" = 10
    "Y24" = 9.5
    "Z24" = 19
    "AA24" = 19
    "AH24" = 17
    "AK24" = 29
    "AO24" = 12
    "AS24" = 201
    "O26" = 1.1
    "P26" = 7
    "G27" = 2.18
    "H27" = 3.1
    "I27" = 3.15
    "J27" = 2.77
    "K27" = 2.05
    "L27" = 3.8
    "N27" = 6.6
    "O27" = 1.36
    "P27" = 2.9
    "Q27" = 2.07
    "R27" = 1.7
    "S27" = 1.42
    "T27" = 2.65
    "X27" = 10.25
    "Z27" = 21
    "AA27" = 18.5
    "AC27" = 6.6
    "AD27" = 6.1
    "AE27" = 14.5
    "AG27" = 8.75
    "AH27" = 16
    "AI27" = 11.25
    "AN27" = 4.05
    "AO27" = 11.5
    "AT27" = 2.65
    "AU27" = 7.1
    "AV27" = 65
    "AW27" = 5.1
    "AX27" = 18
    "AZ27" = 90
    "BA27" = 120
    "G30" = 1.88
    "H30" = 3.75
    "I30" = 3.35
    "K30" = 2.4
    "P30" = 4.75
    "U30" = 1.45
    "V30" = 2.55
    "W30" = 11.5
    "Z30" = 18.5
    "AB30" = 17.5
    "AD30" = 8
    "AE30" = 11.5
    "AJ30" = 50
    "AL30" = 24
    "AM30" = 175
    "AU30" = 6.3
    "AV30" = 37
    "AW30" = 5.9
    "AY30" = 18.5
    "H31" = 3.55
    "I31" = 4.3
    "K31" = 2.18
    "L31" = 4.75
    "M31" = 1.06
    "N31" = 7.4
    "O31" = 1.29
    "P31" = 3.25
    "Q31" = 1.88
    "R31" = 1.85
    "S31" = 1.39
    "T31" = 2.75
    "U31" = 1.83
    "V31" = 1.87
    "W31" = 6.9
    "Z31" = 13.5
    "AB31" = 27
    "AC31" = 7.4
    "AD31" = 7
    "AE31" = 16
    "AF31" = 75
    "AG31" = 12
    "AH31" = 24
    "AJ31" = 75
    "AK31" = 45
    "AL31" = 50
    "AM31" = 600
    "AN31" = 3.55
    "AO31" = 8.5
    "AP31" = 18
    "AT31" = 2.75
    "AU31" = 7.6
    "AY31" = 32
    "BA31" = 200
    "BB31" = 450
    "G32" = 2.25
    "H32" = 3.25
    "I32" = 2.9
    "J32" = 2.87
    "K32" = 2.12
    "O32" = 1.28
    "P32" = 3.35
    "Q32" = 1.85
    "R32" = 1.88
    "S32" = 1.39
    "T32" = 2.77
    "U32" = 1.7
    "V32" = 2.05
    "W32" = 8.25
    "X32" = 11.25
    "Z32" = 23
    "AB32" = 27
    "AD32" = 6.4
    "AG32" = 10
    "AH32" = 16
    "AI32" = 10.5
    "AJ32" = 37
    "AL32" = 29
    "AN32" = 4.25
    "AO32" = 12
    "AT32" = 2.77
    "AU32" = 6.8
    "AV32" = 55
    "AW32" = 4.9
    "AX32" = 15.5
    "AY32" = 22
    "AZ32" = 70
    "BA32" = 100
    "BB32" = 250
    "M35" = 1.03
    "N35" = 12
    "T35" = 2.52
    "AG35" = 15
    "AK35" = 50
    "AL35" = 45
    "AO35" = 8.5
    "AU35" = 7.1
    "AX35" = 29
    "AZ35" = 175
    "G36" = 2.82
    "I36" = 2.55
    "J36" = 3.35
    "K36" = 1.98
    "L36" = 3.2
    "T36" = 2.47
    "Y36" = 9.75
    "AA36" = 23
    "AB36" = 27
    "AG36" = 8.25
    "AH36" = 13.5
    "AJ36" = 30
    "AL36" = 28
    "AO36" = 15.5
    "AP36" = 20
    "AQ36" = 70
    "AR36" = 90
    "AT36" = 2.45
    "AU36" = 6.3
    "AV36" = 50
    "AW36" = 4.5
    "AY36" = 21
    "AZ36" = 65
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Output ("Updated {0} cells" -f $updates.Count)
